$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "quadratic-svm-score" output with the newly computed
# decision-function values for rows 2-8 (previously placeholder 1's),
# matching the regenerated outputs-HGR-r202 run.
$ws.Range("B2").Value = -0.3464487862609742
$ws.Range("B3").Value = -0.27148925610081465
$ws.Range("B4").Value = -0.4381589717418457
$ws.Range("B5").Value = -0.2360612227466663
$ws.Range("B6").Value = -0.23765333749391226
$ws.Range("B7").Value = -0.29300382725716512
$ws.Range("B8").Value = -0.31310196576098104
